$wb = $excel.ActiveWorkbook

# --- Rename the original "Tabelle1" sheet to "Tabelle 1" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Tabelle 1"

# --- Insert two brand-new sheets "Tabelle 2" and "Tabelle 3" right after "Tabelle 1" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Tabelle 2"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Tabelle 3"

$rohdaten = $wb.Worksheets.Item("Rohdaten")

# --- Page setup (paper size / orientation) for the new sheets ---
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

# Tabelle 3 uses different (metric/A4-style) default page margins
$ws3.PageSetup.LeftMargin = 0.70866141732283472 * 72
$ws3.PageSetup.RightMargin = 0.70866141732283472 * 72
$ws3.PageSetup.TopMargin = 0.74803149606299213 * 72
$ws3.PageSetup.BottomMargin = 0.74803149606299213 * 72
$ws3.PageSetup.HeaderMargin = 0.31496062992125984 * 72
$ws3.PageSetup.FooterMargin = 0.31496062992125984 * 72

# --- Print areas: evaluate multiple comma-separated print-area expressions per sheet ---
$ws1.PageSetup.PrintArea = "IF(Rohdaten!`$A`$2<5,'Tabelle 1'!`$A`$1:`$C`$40,'Tabelle 1'!`$A`$1:`$C`$80)"
$ws2.PageSetup.PrintArea = "IF(Rohdaten!`$A`$2<5,'Tabelle 2'!`$A`$1:`$C`$40,'Tabelle 2'!`$A`$1:`$C`$80),'Tabelle 2'!`$J`$17:`$M`$29"
$ws3.PageSetup.PrintArea = "'Tabelle 3'!`$C`$4:`$G`$12,'Tabelle 3'!`$F`$17:`$J`$23"

# --- Selection / active sheet state ---
[void]$rohdaten.Range("A2").Select()

[void]$ws2.Activate()
[void]$ws2.Range("G32").Select()
